$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora): every data row 2-51 changes from "5" to "6"
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "6"

# Column D (Price) updates - only rows whose price changed
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.48"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.70"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.350"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08147"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.612"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.671"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.186"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1349"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1959"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09568"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1049"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001333"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005987"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.398"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.438"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3394"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.158"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3053"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04314"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004269"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001350"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003540"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02774"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05605"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006302"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007689"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1447"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007684"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008089"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3515"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006824"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004001"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"

# Column E (Volume(1h)) updates - only rows whose volume changed
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "8.63%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "18.58%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "6.65%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "8.33%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.40%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.62%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "28.17%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "12.32%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.30%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.79%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11.18%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.09%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.13%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.33%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.23%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.55%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.45%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.15%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.11%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.57%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.48%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.11%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "8.38%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.69%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.93%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "15.18%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.67%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.07%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.63%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.13%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.84%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.25%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "18.75%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.40%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "36.59%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.81%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
